# Update "想去人数" (people interested) counts per gh-pages regeneration (456a3b4)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 495
$ws.Range("F5").Value = 2289
$ws.Range("F8").Value = 75
$ws.Range("F9").Value = 1676
$ws.Range("F10").Value = 1676
$ws.Range("F11").Value = 1383
$ws.Range("F12").Value = 71
$ws.Range("F15").Value = 23
$ws.Range("F16").Value = 672
$ws.Range("F18").Value = 123
$ws.Range("F19").Value = 7403
$ws.Range("F20").Value = 8268
$ws.Range("F21").Value = 54
$ws.Range("F22").Value = 6
$ws.Range("F28").Value = 268
$ws.Range("F31").Value = 4
$ws.Range("F33").Value = 355
$ws.Range("F34").Value = 1477
$ws.Range("F37").Value = 23
$ws.Range("F38").Value = 298
$ws.Range("F40").Value = 762
$ws.Range("F41").Value = 26
$ws.Range("F43").Value = 360
$ws.Range("F44").Value = 259
$ws.Range("F47").Value = 193
$ws.Range("F48").Value = 178
$ws.Range("F49").Value = 18
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 36
$ws.Range("F5").Value = 65
$ws.Range("F13").Value = 10
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 292
$ws.Range("F6").Value = 19
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 36
$ws.Range("F8").Value = 495
$ws.Range("F9").Value = 2289
$ws.Range("F11").Value = 75
$ws.Range("F12").Value = 1676
$ws.Range("F13").Value = 1676
$ws.Range("F14").Value = 71
$ws.Range("F15").Value = 23
$ws.Range("F16").Value = 672
$ws.Range("F19").Value = 65
$ws.Range("F20").Value = 123
$ws.Range("F21").Value = 7403
$ws.Range("F22").Value = 8268
$ws.Range("F23").Value = 54
$ws.Range("F26").Value = 268
$ws.Range("F31").Value = 23
$ws.Range("F33").Value = 298
$ws.Range("F37").Value = 762
$ws.Range("F39").Value = 26
$ws.Range("F40").Value = 10
$ws.Range("F43").Value = 360
$ws.Range("F44").Value = 259
$ws.Range("F47").Value = 193
